$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Motherboard")

# Re-save under the English builtin style name instead of the Dutch
# localized one ("Standaard" -> "Normal"), as Excel does when a workbook
# created under a Dutch locale is opened/saved from an English install
$wb.Styles.Item("Standaard").Delete()

# Replace the placeholder test data in row 3 with real motherboard info
$ws.Range("A3").Value = "X99 Rampage V Extreme"
$ws.Range("C3").Value = "Extended ATX"

# Widen column A to fit the new, longer motherboard name and
# move the active selection from F5 to C5
$ws.Columns.Item(1).ColumnWidth = 22
$ws.Range("C5").Select()
